$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Main" (sheet1) - redesigned toolbar / layout
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Main")

# Content: move "Color info" up to A1, keep B1 "Toolbar", clear A2,
# keep B2 "Diagram", keep A3 "CIE table", clear B3 (stays blank, merged).
$ws1.Range("A1").Value = "Color info"
$ws1.Range("B1").Value = "Toolbar"
$ws1.Range("A2").Value = ""
$ws1.Range("B2").Value = "Diagram"
$ws1.Range("A3").Value = "CIE table"
$ws1.Range("B3").Value = ""

# Merge first, so a uniform border/alignment pass below lands on every
# physical cell of the merged ranges identically.
$ws1.Range("A1:A2").Merge()
$ws1.Range("B2:B3").Merge()

# Borders: thin box around every cell of the table (A1:B3).
$ws1.Range("A1:B3").Borders.LineStyle = 1

# A1/A2 are merged into one vertical header cell - drop the inner seam
# (A1 keeps its top/left/right border, A2 keeps its bottom/left/right).
$ws1.Range("A1").Borders.Item(9).LineStyle = 0
$ws1.Range("A2").Borders.Item(8).LineStyle = 0

# Alignment - every cell centered both ways.
$ws1.Range("A1:B3").HorizontalAlignment = -4108
$ws1.Range("A1:B3").VerticalAlignment = -4108

# Row heights for the redesigned layout.
$ws1.Rows.Item(2).RowHeight = 44.25
$ws1.Rows.Item(3).RowHeight = 209.25

# -----------------------------------------------------------------
# Sheet "Color info" (sheet2) - add borders + spectral distribution row
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Color info")

# New row 5: merged label cell for the spectral distribution chart.
$ws2.Range("A5:C5").Merge()
$ws2.Range("A5").Value = " Spectral distribution"

# Borders: thin box around every used cell, including the existing
# data table (A1:C4) whose B column previously had no border/style at all,
# plus the freshly merged A5:C5 row.
$ws2.Range("A1:C5").Borders.LineStyle = 1

# Alignment - C column (the color swatch + label) + new row stay centered;
# A4:B4 keeps its left alignment, C4 keeps horizontal-center alignment
# (both already set from the source file, preserved, just gain a border).
$ws2.Range("C1:C3").HorizontalAlignment = -4108
$ws2.Range("C1:C3").VerticalAlignment = -4108
$ws2.Range("A5:C5").HorizontalAlignment = -4108
$ws2.Range("A5:C5").VerticalAlignment = -4108

$ws2.Rows.Item(5).RowHeight = 45.75

# -----------------------------------------------------------------
# Selection / active sheet - finish on "Main", cell D3 selected (the
# explicit selection left over on "Color info" gets reset to the sheet's
# top-left cell since it's no longer the one being worked on).
# -----------------------------------------------------------------
$ws2.Range("A1").Select()
$ws1.Select()
$ws1.Range("D3").Select()
